$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.945.30"
$ws.Range("E2").Value = "  -1.57%  "

$ws.Range("D3").Value = "2.505.43"
$ws.Range("E3").Value = "  -4.07%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'577.02"
$ws.Range("E5").Value = "  -2.74%  "

$ws.Range("D6").Value = "'166.45"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'0.521"
$ws.Range("E8").Value = "  -2.07%  "

$ws.Range("D9").Value = "2.503.99"
$ws.Range("E9").Value = "  -4.11%  "

$ws.Range("E10").Value = "  -0.67%  "

$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").Value = "'0.348"
$ws.Range("E12").Value = "  -3.77%  "

$ws.Range("D13").Value = "'5.10"
$ws.Range("E13").Value = "  -2.24%  "

$ws.Range("D14").Value = "'26.31"
$ws.Range("E14").Value = "  -4.57%  "

$ws.Range("D15").Value = "2.962.13"
$ws.Range("E15").Value = "  -4.56%  "

$ws.Range("E16").Value = "  -3.60%  "

$ws.Range("D17").Value = "65.809.93"
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").Value = "2.559.34"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").Value = "'11.21"
$ws.Range("E19").Value = "  -6.42%  "

$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "  -4.21%  "

$ws.Range("D21").Value = "'344.40"
$ws.Range("E21").Value = "  -3.10%  "

$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "  -2.84%  "

$ws.Range("D23").Value = "'4.54"
$ws.Range("E23").Value = "  -2.18%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'1.92"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("D26").Value = "'68.80"
$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  -2.80%  "

$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("D29").Value = "2.644.69"

$ws.Range("D30").Value = "0.0₃0970"
$ws.Range("E30").Value = "  -2.73%  "

$ws.Range("D31").Value = "'8.09"
$ws.Range("E31").Value = "  +2.95%  "

$ws.Range("D32").Value = "'518.00"
$ws.Range("E32").Value = "  -4.57%  "

$ws.Range("E33").Value = "  -3.02%  "

$ws.Range("D34").Value = "'1.80"
$ws.Range("E34").Value = "  -4.63%  "

$ws.Range("E35").Value = "  -3.24%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.06%  "

# Row 37 and 38 swap (ImmutableX <-> Monero) with updated values
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "'157.50"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.44"
$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("D39").Value = "'18.47"
$ws.Range("E39").Value = "  -2.48%  "

$ws.Range("E40").Value = "  +0.72%  "

$ws.Range("D41").Value = "'0.352"
$ws.Range("E41").Value = "  -3.37%  "

$ws.Range("E42").Value = "  -2.57%  "

$ws.Range("D43").Value = "'5.03"
$ws.Range("E43").Value = "  -2.74%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").Value = "'146.50"
$ws.Range("E46").Value = "  -3.22%  "

# Row 47 and 48 swap (BabyDogeCoin <-> ARBITRUM) with updated values
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'0.552"
$ws.Range("E47").Value = "  -4.07%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0276"
$ws.Range("E48").Value = "  -6.55%  "

$ws.Range("D49").Value = "'3.69"
$ws.Range("E49").Value = "  -1.61%  "

$ws.Range("D50").Value = "'1.70"
$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("D51").Value = "'0.0751"
$ws.Range("E51").Value = "  -2.44%  "
